$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Market share_class_min values for ZEV (row 3) for years 2035-2050
# from 0.99 to 1
$ws.Range("T3:W3").Value = 1

# Fix the selection range on the sheet view (was A1:X8, should be A1:X7)
$ws.Range("A1:X7").Select()
